{"js": "// Edit 1: \"He tomado la certificaci\u00f3n de Azure AZ-900, Scrum, y producto owner essentials.\"\n// -> \"He tomado la certificaci\u00f3n de Azure AZ-900, Scrum, y Producto Owner Essentials. \"\n// (capitalize \"Producto\", \"Owner\", \"Essentials\" and add a trailing space after the period)\nconst certResults = context.document.body.search(\n  \", y producto owner essentials.\",\n  { matchCase: true, matchWholeWord: false }\n);\ncertResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < certResults.items.length; i++) {\n  certResults.items[i].insertText(\n    \", y Producto Owner Essentials. \",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n\n// Edit 2: \"Me siento muy seguro en el desarrollo de software...\"\n// -> \"Me siento muy seguraen el desarrollo de software...\"\nconst confidentResults = context.document.body.search(\n  \"Me siento muy seguro en el desarrollo\",\n  { matchCase: true, matchWholeWord: false }\n);\nconfidentResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < confidentResults.items.length; i++) {\n  confidentResults.items[i].insertText(\n    \"Me siento muy seguraen el desarrollo\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Edit 1: \"He tomado la certificaci\u00f3n de Azure AZ-900, Scrum, y producto owner essentials.\"\n# -> \"He tomado la certificaci\u00f3n de Azure AZ-900, Scrum, y Producto Owner Essentials. \"\n# (capitalize \"Producto\", \"Owner\", \"Essentials\" and add a trailing space after the period)\n$find1 = $d.Content.Find\n$found1 = $find1.Execute(\n    \", y producto owner essentials.\",\n    $true,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    \", y Producto Owner Essentials. \",\n    2\n)\n\n# Edit 2: \"Me siento muy seguro en el desarrollo de software...\"\n# -> \"Me siento muy seguraen el desarrollo de software...\"\n$find2 = $d.Content.Find\n$found2 = $find2.Execute(\n    \"Me siento muy seguro en el desarrollo\",\n    $true,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    \"Me siento muy seguraen el desarrollo\",\n    2\n)\n"}
